$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right after the current last data row (30). Inserting
# (rather than just writing into previously-empty rows) makes Excel carry
# down the existing column formatting - e.g. the "email" column's style and
# the "is_active" boolean column's style - exactly like the rest of the
# table, instead of leaving the new rows unstyled.
$ws.Rows(31).Insert()
$ws.Rows(32).Insert()

# Populate row 32 (John Doe) before row 31 (Jane Smith) so the new shared
# strings are appended to the workbook's string table in "John Doe",
# "john.doe@xyz.com", "Jane Smith", "jane.smith@xyz.com" order.
$ws.Cells.Item(32, 1).Value = 110031
$ws.Cells.Item(32, 2).Value = 9317596767
$ws.Cells.Item(32, 3).Value = "John Doe"
$ws.Cells.Item(32, 4).Value = "john.doe@xyz.com"
$ws.Cells.Item(32, 5).Value = 818876431
$ws.Cells.Item(32, 6).Value = "ACT"
$ws.Cells.Item(32, 7).Value = "eng"
$ws.Cells.Item(32, 8).Value = "PWD"
$ws.Cells.Item(32, 9).Value = $true
$ws.Cells.Item(32, 10).Value = "superadmin"
$ws.Cells.Item(32, 11).Value = "now()"
$ws.Cells.Item(32, 12).Value = "now()"

$ws.Cells.Item(31, 1).Value = 110030
$ws.Cells.Item(31, 2).Value = 9317596768
$ws.Cells.Item(31, 3).Value = "Jane Smith"
$ws.Cells.Item(31, 4).Value = "jane.smith@xyz.com"
$ws.Cells.Item(31, 5).Value = 818876432
$ws.Cells.Item(31, 6).Value = "ACT"
$ws.Cells.Item(31, 7).Value = "eng"
$ws.Cells.Item(31, 8).Value = "PWD"
$ws.Cells.Item(31, 9).Value = $true
$ws.Cells.Item(31, 10).Value = "superadmin"
$ws.Cells.Item(31, 11).Value = "now()"
$ws.Cells.Item(31, 12).Value = "now()"

# Leave the selection where the editor last left it.
$ws.Range("F30").Select()
